# Applies the "updated results, added discussion from rm with my edits" revision:
#   * Each Heading1/Heading2 section bookmark (introduction, methods,
#     general-approach, stressor-components, response-components,
#     stream-condition-index, application, results, discussion) is widened
#     so it spans the heading text itself instead of sitting collapsed in
#     front of it.
#   * An (empty) comments part is introduced for the document.
#
# Note: hyperlink/footer relationship-id renumbering and the numbering.xml
# abstractNum bookkeeping seen in the source diff are save-time side
# effects of the authoring tool's relationship-id allocator, not user
# edits; they carry no addressable semantics through the Word object
# model, so they are intentionally left alone here.

$d = $word.ActiveDocument

function Expand-SectionBookmark($name) {
    $bm = $d.Bookmarks.Item($name)
    $para = $bm.Range.Paragraphs(1)
    $paraStart = $para.Range.Start
    $paraEnd = $para.Range.End

    # Heading paragraphs here hold a single run with the heading text,
    # followed by the paragraph mark; trim that trailing mark off so the
    # bookmark wraps just the visible text (matches target structure:
    # bookmarkStart, run, bookmarkEnd).
    $headingRange = $d.Range($paraStart, $paraEnd - 1)

    $bm.Delete()
    $d.Bookmarks.Add($name, $headingRange) | Out-Null
}

$sectionNames = @(
    "introduction",
    "methods",
    "general-approach",
    "stressor-components",
    "response-components",
    "stream-condition-index",
    "application",
    "results",
    "discussion"
)

foreach ($name in $sectionNames) {
    Expand-SectionBookmark $name
}

# Introduce the (empty) comments part/relationship for this document.
$commentRange = $d.Range(0, 0)
$tempComment = $d.Comments.Add($commentRange, " ")
$tempComment.Delete()

Write-Output "done"
